$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

# New rows 205-227: mmWave presence-sensor log entries for 2026-01-28.
# Column A holds an ISO-formatted date string ("2026-01-28"). Excel
# would normally auto-convert that to a date serial on assignment, but
# every existing row in this sheet stores the date as plain text, so we
# temporarily force Text format before writing it, then clear the
# formatting again afterwards so the new cells end up styled exactly
# like the rest of the sheet (default/general style, text content).
$dateCol = $ws.Range("A205:A227")
$dateCol.NumberFormat = "@"

$ws.Cells.Item(205, 1).Value = "2026-01-28"
$ws.Cells.Item(205, 2).Value = "17:51:49"
$ws.Cells.Item(205, 3).Value = "17:00"
$ws.Cells.Item(205, 4).Value = "Living Room"
$ws.Cells.Item(205, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(205, 6).Value = "Inactive"

$ws.Cells.Item(206, 1).Value = "2026-01-28"
$ws.Cells.Item(206, 2).Value = "17:51:49"
$ws.Cells.Item(206, 3).Value = "17:00"
$ws.Cells.Item(206, 4).Value = "Living Room"
$ws.Cells.Item(206, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(206, 6).Value = "Inactive"

$ws.Cells.Item(207, 1).Value = "2026-01-28"
$ws.Cells.Item(207, 2).Value = "17:51:50"
$ws.Cells.Item(207, 3).Value = "17:00"
$ws.Cells.Item(207, 4).Value = "Living Room"
$ws.Cells.Item(207, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(207, 6).Value = "Inactive"

$ws.Cells.Item(208, 1).Value = "2026-01-28"
$ws.Cells.Item(208, 2).Value = "17:51:51"
$ws.Cells.Item(208, 3).Value = "17:00"
$ws.Cells.Item(208, 4).Value = "Living Room"
$ws.Cells.Item(208, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(208, 6).Value = "Inactive"

$ws.Cells.Item(209, 1).Value = "2026-01-28"
$ws.Cells.Item(209, 2).Value = "17:51:55"
$ws.Cells.Item(209, 3).Value = "17:00"
$ws.Cells.Item(209, 4).Value = "Living Room"
$ws.Cells.Item(209, 5).Value = "PRESENCE"
$ws.Cells.Item(209, 6).Value = "Active"

$ws.Cells.Item(210, 1).Value = "2026-01-28"
$ws.Cells.Item(210, 2).Value = "17:51:58"
$ws.Cells.Item(210, 3).Value = "17:00"
$ws.Cells.Item(210, 4).Value = "Living Room"
$ws.Cells.Item(210, 5).Value = "PRESENCE"
$ws.Cells.Item(210, 6).Value = "Active"

$ws.Cells.Item(211, 1).Value = "2026-01-28"
$ws.Cells.Item(211, 2).Value = "17:52:01"
$ws.Cells.Item(211, 3).Value = "17:00"
$ws.Cells.Item(211, 4).Value = "Living Room"
$ws.Cells.Item(211, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(211, 6).Value = "Inactive"

$ws.Cells.Item(212, 1).Value = "2026-01-28"
$ws.Cells.Item(212, 2).Value = "17:52:04"
$ws.Cells.Item(212, 3).Value = "17:00"
$ws.Cells.Item(212, 4).Value = "Living Room"
$ws.Cells.Item(212, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(212, 6).Value = "Inactive"

$ws.Cells.Item(213, 1).Value = "2026-01-28"
$ws.Cells.Item(213, 2).Value = "17:52:07"
$ws.Cells.Item(213, 3).Value = "17:00"
$ws.Cells.Item(213, 4).Value = "Living Room"
$ws.Cells.Item(213, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(213, 6).Value = "Inactive"

$ws.Cells.Item(214, 1).Value = "2026-01-28"
$ws.Cells.Item(214, 2).Value = "17:52:10"
$ws.Cells.Item(214, 3).Value = "17:00"
$ws.Cells.Item(214, 4).Value = "Living Room"
$ws.Cells.Item(214, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(214, 6).Value = "Inactive"

$ws.Cells.Item(215, 1).Value = "2026-01-28"
$ws.Cells.Item(215, 2).Value = "17:52:13"
$ws.Cells.Item(215, 3).Value = "17:00"
$ws.Cells.Item(215, 4).Value = "Living Room"
$ws.Cells.Item(215, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(215, 6).Value = "Inactive"

$ws.Cells.Item(216, 1).Value = "2026-01-28"
$ws.Cells.Item(216, 2).Value = "17:52:16"
$ws.Cells.Item(216, 3).Value = "17:00"
$ws.Cells.Item(216, 4).Value = "Living Room"
$ws.Cells.Item(216, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(216, 6).Value = "Inactive"

$ws.Cells.Item(217, 1).Value = "2026-01-28"
$ws.Cells.Item(217, 2).Value = "17:52:19"
$ws.Cells.Item(217, 3).Value = "17:00"
$ws.Cells.Item(217, 4).Value = "Living Room"
$ws.Cells.Item(217, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(217, 6).Value = "Inactive"

$ws.Cells.Item(218, 1).Value = "2026-01-28"
$ws.Cells.Item(218, 2).Value = "17:52:22"
$ws.Cells.Item(218, 3).Value = "17:00"
$ws.Cells.Item(218, 4).Value = "Living Room"
$ws.Cells.Item(218, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(218, 6).Value = "Inactive"

$ws.Cells.Item(219, 1).Value = "2026-01-28"
$ws.Cells.Item(219, 2).Value = "17:52:25"
$ws.Cells.Item(219, 3).Value = "17:00"
$ws.Cells.Item(219, 4).Value = "Living Room"
$ws.Cells.Item(219, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(219, 6).Value = "Inactive"

$ws.Cells.Item(220, 1).Value = "2026-01-28"
$ws.Cells.Item(220, 2).Value = "17:52:28"
$ws.Cells.Item(220, 3).Value = "17:00"
$ws.Cells.Item(220, 4).Value = "Living Room"
$ws.Cells.Item(220, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(220, 6).Value = "Inactive"

$ws.Cells.Item(221, 1).Value = "2026-01-28"
$ws.Cells.Item(221, 2).Value = "17:52:31"
$ws.Cells.Item(221, 3).Value = "17:00"
$ws.Cells.Item(221, 4).Value = "Living Room"
$ws.Cells.Item(221, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(221, 6).Value = "Inactive"

$ws.Cells.Item(222, 1).Value = "2026-01-28"
$ws.Cells.Item(222, 2).Value = "17:52:34"
$ws.Cells.Item(222, 3).Value = "17:00"
$ws.Cells.Item(222, 4).Value = "Living Room"
$ws.Cells.Item(222, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(222, 6).Value = "Inactive"

$ws.Cells.Item(223, 1).Value = "2026-01-28"
$ws.Cells.Item(223, 2).Value = "17:52:37"
$ws.Cells.Item(223, 3).Value = "17:00"
$ws.Cells.Item(223, 4).Value = "Living Room"
$ws.Cells.Item(223, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(223, 6).Value = "Inactive"

$ws.Cells.Item(224, 1).Value = "2026-01-28"
$ws.Cells.Item(224, 2).Value = "17:52:40"
$ws.Cells.Item(224, 3).Value = "17:00"
$ws.Cells.Item(224, 4).Value = "Living Room"
$ws.Cells.Item(224, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(224, 6).Value = "Inactive"

$ws.Cells.Item(225, 1).Value = "2026-01-28"
$ws.Cells.Item(225, 2).Value = "17:52:43"
$ws.Cells.Item(225, 3).Value = "17:00"
$ws.Cells.Item(225, 4).Value = "Living Room"
$ws.Cells.Item(225, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(225, 6).Value = "Inactive"

$ws.Cells.Item(226, 1).Value = "2026-01-28"
$ws.Cells.Item(226, 2).Value = "17:52:46"
$ws.Cells.Item(226, 3).Value = "17:00"
$ws.Cells.Item(226, 4).Value = "Living Room"
$ws.Cells.Item(226, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(226, 6).Value = "Inactive"

$ws.Cells.Item(227, 1).Value = "2026-01-28"
$ws.Cells.Item(227, 2).Value = "17:52:49"
$ws.Cells.Item(227, 3).Value = "17:00"
$ws.Cells.Item(227, 4).Value = "Living Room"
$ws.Cells.Item(227, 5).Value = "NO_PRESENCE"
$ws.Cells.Item(227, 6).Value = "Inactive"

$dateCol.ClearFormats()